# Generate Report for Handback
# Reorders the status rows (the two files that were "Handed back" move to
# the top of each sheet) and fills in the handback columns (Latest Target
# File / Latest Handback File / Latest Handback DateTime) for those rows.

$wb = $excel.ActiveWorkbook

$mdUrl = @{
  "05976f76" = "https://github.com/OpenLocalizationTest/oltest/blob/2f07166738f831f7e0dcfc33c45c9407a7eda7e1/e2e/05976f76-c427-4154-b354-7b80eedb385f.md"
  "6b9de217" = "https://github.com/OpenLocalizationTest/oltest/blob/d544eea60da2b4c32f5ebd742d344626d640ee78/e2e/6b9de217-6b4d-429d-880d-6db2f2fd5d79.md"
  "6cf0ffb6" = "https://github.com/OpenLocalizationTest/oltest/blob/2f07166738f831f7e0dcfc33c45c9407a7eda7e1/e2e/6cf0ffb6-b5f7-4b18-a18e-29508156f435.md"
  "81191734" = "https://github.com/OpenLocalizationTest/oltest/blob/d544eea60da2b4c32f5ebd742d344626d640ee78/e2e/81191734-0aaf-40e7-b3f4-2a68a2dbad53.md"
}
$mdName = @{
  "05976f76" = "05976f76-c427-4154-b354-7b80eedb385f.md"
  "6b9de217" = "6b9de217-6b4d-429d-880d-6db2f2fd5d79.md"
  "6cf0ffb6" = "6cf0ffb6-b5f7-4b18-a18e-29508156f435.md"
  "81191734" = "81191734-0aaf-40e7-b3f4-2a68a2dbad53.md"
}
$zhXlfUrl = @{
  "05976f76" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7e7ebcf3060337814b511aba6d74774d0610803/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/05976f76-c427-4154-b354-7b80eedb385f.a79645bbac857611d2e28ccb52d1601f78bc72f1.zh-cn.xlf"
  "6b9de217" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7e7ebcf3060337814b511aba6d74774d0610803/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/6b9de217-6b4d-429d-880d-6db2f2fd5d79.7a0af47e61cb8dd5c48a95f258828a6e50bcca54.zh-cn.xlf"
  "6cf0ffb6" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7e7ebcf3060337814b511aba6d74774d0610803/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/6cf0ffb6-b5f7-4b18-a18e-29508156f435.6f3b041ab5d91f9d71dbafb9fa9676707d9f24e4.zh-cn.xlf"
  "81191734" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7e7ebcf3060337814b511aba6d74774d0610803/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/81191734-0aaf-40e7-b3f4-2a68a2dbad53.9a9425152b91a59b900a463a384f59f70e7c7653.zh-cn.xlf"
}
$zhXlfName = @{
  "05976f76" = "05976f76-c427-4154-b354-7b80eedb385f.a79645bbac857611d2e28ccb52d1601f78bc72f1.zh-cn.xlf"
  "6b9de217" = "6b9de217-6b4d-429d-880d-6db2f2fd5d79.7a0af47e61cb8dd5c48a95f258828a6e50bcca54.zh-cn.xlf"
  "6cf0ffb6" = "6cf0ffb6-b5f7-4b18-a18e-29508156f435.6f3b041ab5d91f9d71dbafb9fa9676707d9f24e4.zh-cn.xlf"
  "81191734" = "81191734-0aaf-40e7-b3f4-2a68a2dbad53.9a9425152b91a59b900a463a384f59f70e7c7653.zh-cn.xlf"
}
$deXlfUrl = @{
  "05976f76" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ee8beb0e20263662ae917a3b7041da7f3133047/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/05976f76-c427-4154-b354-7b80eedb385f.a79645bbac857611d2e28ccb52d1601f78bc72f1.de-de.xlf"
  "6b9de217" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ee8beb0e20263662ae917a3b7041da7f3133047/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/6b9de217-6b4d-429d-880d-6db2f2fd5d79.7a0af47e61cb8dd5c48a95f258828a6e50bcca54.de-de.xlf"
  "6cf0ffb6" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ee8beb0e20263662ae917a3b7041da7f3133047/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/6cf0ffb6-b5f7-4b18-a18e-29508156f435.6f3b041ab5d91f9d71dbafb9fa9676707d9f24e4.de-de.xlf"
  "81191734" = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7ee8beb0e20263662ae917a3b7041da7f3133047/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/81191734-0aaf-40e7-b3f4-2a68a2dbad53.9a9425152b91a59b900a463a384f59f70e7c7653.de-de.xlf"
}
$deXlfName = @{
  "05976f76" = "05976f76-c427-4154-b354-7b80eedb385f.a79645bbac857611d2e28ccb52d1601f78bc72f1.de-de.xlf"
  "6b9de217" = "6b9de217-6b4d-429d-880d-6db2f2fd5d79.7a0af47e61cb8dd5c48a95f258828a6e50bcca54.de-de.xlf"
  "6cf0ffb6" = "6cf0ffb6-b5f7-4b18-a18e-29508156f435.6f3b041ab5d91f9d71dbafb9fa9676707d9f24e4.de-de.xlf"
  "81191734" = "81191734-0aaf-40e7-b3f4-2a68a2dbad53.9a9425152b91a59b900a463a384f59f70e7c7653.de-de.xlf"
}

$handedBack = "Handed back: in sync with en-US"
$readyFor = "Ready for handoff"

# New row order (top -> bottom) for every sheet: the two "handed back" files
# first, then the two still "ready for handoff".
$order = @("6b9de217", "81191734", "05976f76", "6cf0ffb6")

function Style-Hyperlink($rng) {
  $rng.Font.Underline = $true
  $rng.Font.Color = 15570276
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A1").Hyperlinks.Delete()

for ($i = 0; $i -lt 4; $i++) {
  $key = $order[$i]
  $row = $i + 2
  $status = $readyFor
  if ($i -lt 2) { $status = $handedBack }

  $ws1.Cells.Item($row, 1).Value = $mdName[$key]
  $ws1.Cells.Item($row, 2).Value = $status
  $ws1.Cells.Item($row, 3).Value = $status
  $ws1.Cells.Item($row, 4).Value = "2016-17-13 10:17:12"

  $ws1.Hyperlinks.Add($ws1.Cells.Item($row, 1), $mdUrl[$key], "", "", $mdName[$key])
}

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A1").Hyperlinks.Delete()

for ($i = 0; $i -lt 4; $i++) {
  $key = $order[$i]
  $row = $i + 2
  $handedBackRow = ($i -lt 2)
  $status = $readyFor
  if ($handedBackRow) { $status = $handedBack }

  $ws2.Cells.Item($row, 1).Value = $mdName[$key]
  $ws2.Cells.Item($row, 2).Value = ".md"
  $ws2.Cells.Item($row, 3).Value = $status
  $ws2.Cells.Item($row, 4).Value = $zhXlfName[$key]
  $ws2.Cells.Item($row, 5).Value = "2016-03-13 10:17:08"
  $ws2.Cells.Item($row, 8).Value = "0001-01-01 00:00:00"
  $ws2.Cells.Item($row, 9).Value = "Include"

  $ws2.Hyperlinks.Add($ws2.Cells.Item($row, 1), $mdUrl[$key], "", "", $mdName[$key])
  $ws2.Hyperlinks.Add($ws2.Cells.Item($row, 2), $mdUrl[$key], "", "", ".md")
  $ws2.Hyperlinks.Add($ws2.Cells.Item($row, 4), $zhXlfUrl[$key], "", "", $zhXlfName[$key])

  if ($handedBackRow) {
    $ws2.Cells.Item($row, 6).Value = $mdName[$key]
    $ws2.Cells.Item($row, 7).Value = $zhXlfName[$key]
    $ws2.Cells.Item($row, 8).Value = "2016-03-13 10:19:17"

    Style-Hyperlink($ws2.Cells.Item($row, 6))
    Style-Hyperlink($ws2.Cells.Item($row, 7))

    $ws2.Hyperlinks.Add($ws2.Cells.Item($row, 6), $mdUrl[$key], "", "", $mdName[$key])
    $ws2.Hyperlinks.Add($ws2.Cells.Item($row, 7), $zhXlfUrl[$key], "", "", $zhXlfName[$key])
  }
}

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A1").Hyperlinks.Delete()

for ($i = 0; $i -lt 4; $i++) {
  $key = $order[$i]
  $row = $i + 2
  $handedBackRow = ($i -lt 2)
  $status = $readyFor
  if ($handedBackRow) { $status = $handedBack }

  $ws3.Cells.Item($row, 1).Value = $mdName[$key]
  $ws3.Cells.Item($row, 2).Value = ".md"
  $ws3.Cells.Item($row, 3).Value = $status
  $ws3.Cells.Item($row, 4).Value = $deXlfName[$key]
  $ws3.Cells.Item($row, 5).Value = "2016-03-13 10:17:12"
  $ws3.Cells.Item($row, 8).Value = "0001-01-01 00:00:00"
  $ws3.Cells.Item($row, 9).Value = "Include"

  $ws3.Hyperlinks.Add($ws3.Cells.Item($row, 1), $mdUrl[$key], "", "", $mdName[$key])
  $ws3.Hyperlinks.Add($ws3.Cells.Item($row, 2), $mdUrl[$key], "", "", ".md")
  $ws3.Hyperlinks.Add($ws3.Cells.Item($row, 4), $deXlfUrl[$key], "", "", $deXlfName[$key])

  if ($handedBackRow) {
    $ws3.Cells.Item($row, 6).Value = $mdName[$key]
    $ws3.Cells.Item($row, 7).Value = $deXlfName[$key]
    $ws3.Cells.Item($row, 8).Value = "2016-03-13 10:19:23"

    Style-Hyperlink($ws3.Cells.Item($row, 6))
    Style-Hyperlink($ws3.Cells.Item($row, 7))

    $ws3.Hyperlinks.Add($ws3.Cells.Item($row, 6), $mdUrl[$key], "", "", $mdName[$key])
    $ws3.Hyperlinks.Add($ws3.Cells.Item($row, 7), $deXlfUrl[$key], "", "", $deXlfName[$key])
  }
}
